$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "mostFrom" shared text value for all data rows (Singapore -> Selangor)
$ws.Range("D2:D6").Value = "Selangor"

# Update avgPrice (column B) and avgSold (column C) values
$ws.Range("B2").Value = 20.9596551724138
$ws.Range("C2").Value = 288.51724137931

$ws.Range("B3").Value = 16.1383870967742
$ws.Range("C3").Value = 1122.77419354839

$ws.Range("B4").Value = 11.84
$ws.Range("C4").Value = 3523.38461538462

$ws.Range("B5").Value = 17.6
$ws.Range("C5").Value = 396.875

$ws.Range("B6").Value = 6.77914285714286
$ws.Range("C6").Value = 1862.54285714286
